$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.697.08'
$ws.Range('E2').Value = '  +6.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.450.30'
$ws.Range('E3').Value = '  +8.77%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '477.79'
$ws.Range('E5').Value = '  +11.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.63'
$ws.Range('E6').Value = '  +22.20%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.501'
$ws.Range('E8').Value = '  +12.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.450.52'
$ws.Range('E9').Value = '  +8.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0958'
$ws.Range('E10').Value = '  +16.93%  '
$ws.Range('E11').Value = '  +7.49%  '
$ws.Range('E12').Value = '  +11.23%  '
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.845.14'
$ws.Range('E14').Value = '  +7.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '54.905.41'
$ws.Range('E15').Value = '  +6.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.42'
$ws.Range('E16').Value = '  +14.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  +21.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.446.69'
$ws.Range('E18').Value = '  +8.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.34'
$ws.Range('E19').Value = '  +14.67%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '311.79'
$ws.Range('E20').Value = '  +8.93%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.80'
$ws.Range('E21').Value = '  +18.14%  '
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.63'
$ws.Range('E23').Value = '  +17.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.04'
$ws.Range('E24').Value = '  +9.74%  '
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('E26').Value = '  +13.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.159'
$ws.Range('E27').Value = '  +20.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.543.68'
$ws.Range('E28').Value = '  +8.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.31'
$ws.Range('E29').Value = '  +12.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0763'
$ws.Range('E30').Value = '  +26.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.996'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.40'
$ws.Range('E32').Value = '  +3.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.85'
$ws.Range('E33').Value = '  +11.35%  '
$ws.Range('E34').Value = '  +15.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.14'
$ws.Range('E35').Value = '  +14.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.11'
$ws.Range('E36').Value = '  +19.08%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.58'
$ws.Range('E37').Value = '  +12.03%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.841'
$ws.Range('E38').Value = '  +11.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.40'
$ws.Range('E39').Value = '  +6.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.996'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.600'
$ws.Range('E41').Value = '  +10.45%  '
$ws.Range('E42').Value = '  +14.51%  '
$ws.Range('E43').Value = '  +13.69%  '
$ws.Range('E44').Value = '  +18.66%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '254.76'
$ws.Range('E46').Value = '  +37.64%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.63'
$ws.Range('E47').Value = '  +22.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0887'
$ws.Range('E48').Value = '  +14.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0220'
$ws.Range('E49').Value = '  +13.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.911.46'
$ws.Range('E50').Value = '  +4.85%  '
$ws.Range('E51').Value = '  +14.20%  '
